# Daily TGP (terminal gate pricing) refresh: shift the rolling "Effective Date"
# window forward a day per state block and update the corresponding Diesel/
# ULP/PULP/e10 prices (cents/litre) to the latest published figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A8").Value = 46035
$ws.Range("D8").Value = 152.91999999999999
$ws.Range("E8").Value = 149.5
$ws.Range("F8").Value = 159.5
$ws.Range("G8").Value = 149.52000000000001
$ws.Range("A9").Value = 46035
$ws.Range("D9").Value = 152.91999999999999
$ws.Range("E9").Value = 149.5
$ws.Range("F9").Value = 159.5
$ws.Range("G9").Value = 149.52000000000001
$ws.Range("A10").Value = 46035
$ws.Range("D10").Value = 154.04
$ws.Range("E10").Value = 151.15
$ws.Range("F10").Value = 161.15
$ws.Range("G10").Value = 151.55000000000001
$ws.Range("A11").Value = 46032
$ws.Range("D11").Value = 152.83000000000001
$ws.Range("E11").Value = 149.25
$ws.Range("F11").Value = 159.25
$ws.Range("G11").Value = 149.27000000000001
$ws.Range("A12").Value = 46032
$ws.Range("D12").Value = 152.83000000000001
$ws.Range("E12").Value = 149.25
$ws.Range("F12").Value = 159.25
$ws.Range("G12").Value = 149.27000000000001
$ws.Range("A13").Value = 46032
$ws.Range("D13").Value = 153.76
$ws.Range("E13").Value = 150.69999999999999
$ws.Range("F13").Value = 160.69999999999999
$ws.Range("G13").Value = 151.1
$ws.Range("A17").Value = 46035
$ws.Range("D17").Value = 157.32
$ws.Range("E17").Value = 153.83000000000001
$ws.Range("F17").Value = 163.83000000000001
$ws.Range("A18").Value = 46032
$ws.Range("D18").Value = 157.26
$ws.Range("E18").Value = 153.72
$ws.Range("F18").Value = 163.72
$ws.Range("A22").Value = 46035
$ws.Range("D22").Value = 153.88
$ws.Range("E22").Value = 150.63999999999999
$ws.Range("F22").Value = 160.24
$ws.Range("G22").Value = 151.71
$ws.Range("A23").Value = 46035
$ws.Range("D23").Value = 158.46
$ws.Range("E23").Value = 156.49
$ws.Range("F23").Value = 166.49
$ws.Range("A24").Value = 46035
$ws.Range("D24").Value = 158.62
$ws.Range("E24").Value = 157.12
$ws.Range("F24").Value = 167.12
$ws.Range("A25").Value = 46035
$ws.Range("D25").Value = 158.61000000000001
$ws.Range("E25").Value = 156.63999999999999
$ws.Range("F25").Value = 166.64
$ws.Range("G25").Value = 156.77000000000001
$ws.Range("A26").Value = 46035
$ws.Range("D26").Value = 158.24
$ws.Range("E26").Value = 158.25
$ws.Range("F26").Value = 168.25
$ws.Range("A27").Value = 46032
$ws.Range("D27").Value = 153.69
$ws.Range("E27").Value = 150.5
$ws.Range("F27").Value = 160.1
$ws.Range("G27").Value = 151.57
$ws.Range("A28").Value = 46032
$ws.Range("D28").Value = 158.28
$ws.Range("E28").Value = 156.15
$ws.Range("F28").Value = 166.15
$ws.Range("A29").Value = 46032
$ws.Range("D29").Value = 158.44
$ws.Range("E29").Value = 156.78
$ws.Range("F29").Value = 166.78
$ws.Range("A30").Value = 46032
$ws.Range("D30").Value = 158.43
$ws.Range("E30").Value = 156.29
$ws.Range("F30").Value = 166.29
$ws.Range("G30").Value = 156.41999999999999
$ws.Range("A31").Value = 46032
$ws.Range("D31").Value = 158.07
$ws.Range("E31").Value = 157.9
$ws.Range("F31").Value = 167.9
$ws.Range("A35").Value = 46035
$ws.Range("D35").Value = 152.49
$ws.Range("E35").Value = 148.13999999999999
$ws.Range("F35").Value = 157.13999999999999
$ws.Range("A36").Value = 46032
$ws.Range("D36").Value = 152.09
$ws.Range("E36").Value = 148.01
$ws.Range("F36").Value = 157.01
$ws.Range("A40").Value = 46035
$ws.Range("D40").Value = 157.97999999999999
$ws.Range("E40").Value = 156.43
$ws.Range("F40").Value = 166.43
$ws.Range("A41").Value = 46035
$ws.Range("D41").Value = 157.71
$ws.Range("E41").Value = 156.85
$ws.Range("F41").Value = 166.85
$ws.Range("A42").Value = 46032
$ws.Range("D42").Value = 157.91
$ws.Range("E42").Value = 156.32
$ws.Range("F42").Value = 166.32
$ws.Range("A43").Value = 46032
$ws.Range("D43").Value = 157.63
$ws.Range("E43").Value = 156.74
$ws.Range("F43").Value = 166.74
$ws.Range("A47").Value = 46035
$ws.Range("D47").Value = 152.56
$ws.Range("E47").Value = 149.85
$ws.Range("F47").Value = 159.85
$ws.Range("A48").Value = 46035
$ws.Range("D48").Value = 152.19999999999999
$ws.Range("E48").Value = 149.79
$ws.Range("F48").Value = 159.79
$ws.Range("A49").Value = 46032
$ws.Range("D49").Value = 152.62
$ws.Range("E49").Value = 149.93
$ws.Range("F49").Value = 159.93
$ws.Range("A50").Value = 46032
$ws.Range("D50").Value = 152.26
$ws.Range("E50").Value = 149.87
$ws.Range("F50").Value = 159.87
$ws.Range("A54").Value = 46035
$ws.Range("D54").Value = 167.21
$ws.Range("E54").Value = 163.6
$ws.Range("F54").Value = 173.6
$ws.Range("A55").Value = 46035
$ws.Range("D55").Value = 160.38
$ws.Range("E55").Value = 162.68
$ws.Range("F55").Value = 172.68
$ws.Range("A56").Value = 46035
$ws.Range("D56").Value = 156.72999999999999
$ws.Range("A57").Value = 46035
$ws.Range("D57").Value = 157.38999999999999
$ws.Range("E57").Value = 157.1
$ws.Range("A58").Value = 46035
$ws.Range("D58").Value = 153.15
$ws.Range("E58").Value = 153
$ws.Range("F58").Value = 163
$ws.Range("A59").Value = 46035
$ws.Range("D59").Value = 159.66999999999999
$ws.Range("E59").Value = 161.99
$ws.Range("A60").Value = 46032
$ws.Range("D60").Value = 167.15
$ws.Range("E60").Value = 163.43
$ws.Range("F60").Value = 173.43
$ws.Range("A61").Value = 46032
$ws.Range("D61").Value = 160.29
$ws.Range("E61").Value = 162.57
$ws.Range("F61").Value = 172.57
$ws.Range("A62").Value = 46032
$ws.Range("D62").Value = 156.77000000000001
$ws.Range("A63").Value = 46032
$ws.Range("D63").Value = 157.44
$ws.Range("E63").Value = 156.99
$ws.Range("A64").Value = 46032
$ws.Range("D64").Value = 153.21
$ws.Range("E64").Value = 152.88999999999999
$ws.Range("F64").Value = 162.88999999999999
$ws.Range("A65").Value = 46032
$ws.Range("D65").Value = 159.62
$ws.Range("E65").Value = 161.85
